# Insert a new data row at row 627 (pushing the existing rows 627-642 down
# to 628-643) and populate it with the new "Perejil" price observation,
# matching the weekly update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("627:627").Insert()

$ws.Cells.Item(627, 1).Value = 9
$ws.Cells.Item(627, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(627, 3).Value = "Metropolitana"
$ws.Cells.Item(627, 4).Value = 45239
$ws.Cells.Item(627, 5).Value = 13
$ws.Cells.Item(627, 6).Value = 100112044
$ws.Cells.Item(627, 7).Value = "Perejil"
$ws.Cells.Item(627, 8).Value = "Sin especificar"
$ws.Cells.Item(627, 9).Value = "Primera"
$ws.Cells.Item(627, 10).Value = 70
$ws.Cells.Item(627, 11).Value = 10000
$ws.Cells.Item(627, 12).Value = 12000
$ws.Cells.Item(627, 13).Value = 11000
$ws.Cells.Item(627, 14).Value = "`$/docena de atados"
$ws.Cells.Item(627, 15).Value = "Región Metropolitana"
$ws.Cells.Item(627, 16).Value = 3667
$ws.Cells.Item(627, 17).Value = 3
$ws.Cells.Item(627, 18).Value = "Hortaliza"
